$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 36
$ws.Range("F6").Value = 9139
$ws.Range("F8").Value = 249
$ws.Range("F9").Value = 7177
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 5565
$ws.Range("F14").Value = 6369
$ws.Range("F15").Value = 1101
$ws.Range("F16").Value = 432
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 598
$ws.Range("F20").Value = 280
$ws.Range("F22").Value = 161
$ws.Range("F23").Value = 104
$ws.Range("F24").Value = 10286
$ws.Range("F25").Value = 1953
$ws.Range("F26").Value = 2164
$ws.Range("F27").Value = 47
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 173
$ws.Range("F33").Value = 22
$ws.Range("F35").Value = 2125
$ws.Range("F38").Value = 0
$ws.Range("F39").Value = 5304
$ws.Range("F40").Value = 1207
$ws.Range("F41").Value = 697
$ws.Range("F43").Value = 168
$ws.Range("F45").Value = 0
$ws.Range("F46").Value = 995
$ws.Range("F47").Value = 0

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 6
$ws.Range("F7").Value = 238
$ws.Range("F8").Value = 8
$ws.Range("F9").Value = 45
$ws.Range("F10").Value = 196
$ws.Range("F12").Value = 9
$ws.Range("F13").Value = 2
$ws.Range("F15").Value = 99
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 911
$ws.Range("F20").Value = 0
$ws.Range("F22").Value = 3

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 40
$ws.Range("F3").Value = 36
$ws.Range("F4").Value = 78
$ws.Range("F7").Value = 1167
$ws.Range("F8").Value = 9139
$ws.Range("F10").Value = 249
$ws.Range("F11").Value = 7177
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 6
$ws.Range("F15").Value = 3
$ws.Range("F19").Value = 6369
$ws.Range("F21").Value = 432
$ws.Range("F22").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 280
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 161
$ws.Range("F27").Value = 196
$ws.Range("F28").Value = 10286
$ws.Range("F29").Value = 1953
$ws.Range("F30").Value = 2164
$ws.Range("F35").Value = 22
$ws.Range("F37").Value = 2126
$ws.Range("F38").Value = 313
$ws.Range("F39").Value = 1422
$ws.Range("F41").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("F43").Value = 122
$ws.Range("F46").Value = 0
$ws.Range("F47").Value = 995
$ws.Range("F48").Value = 1391
$ws.Range("F50").Value = 1099
